# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held the literal string "5-27-2013-14" for every
# data row; correct it to the ISO-style "2014-05-27".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value -eq "5-27-2013-14") {
        $cell.Value = "2014-05-27"
    }
}
